$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the resource catalog table with the actual data
$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "Monturas"

$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Lentes"

$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Caja registradora"

$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Datáfono"

$ws.Range("A10").Value = 5
$ws.Range("B10").Value = "Laboratorio"

$ws.Range("A11").Value = 6
$ws.Range("B11").Value = "Consultorio"

$ws.Range("A12").Value = 7
$ws.Range("B12").Value = "Sistema informático"

# The 13th row was left over from the original template and is no longer needed
$ws.Rows.Item(13).Delete()

# Update selection to where the user left off editing
[void]$ws.Range("B22").Select()
